$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows for years 2004-2009 (rows 2-7), shifting everything up.
$ws.Range("A2:F7").Delete(-4162)

# After the shift, years 2010-2020 now occupy rows 2-12.
# Update column B values (more precise figures) for 2016-2020 (rows 8-12).
$ws.Range("B8").Value = 2339.43175840725
$ws.Range("B9").Value = 2059.94062519025
$ws.Range("B10").Value = 1957.74758333571
$ws.Range("B11").Value = 2062.94463841108
$ws.Range("B12").Value = 2239.7562185529

# Add new row 13 for 2021.
$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value = 2098.48763771276
$ws.Range("C13").Value = 8195.700000000001
$ws.Range("D13").Value = 6868
$ws.Range("E13").Value = 28310.5
$ws.Range("F13").Value = 29638.2

# Add new row 14 for 2022, with only the F column populated.
$ws.Range("A14").Value = "2022年"
$ws.Range("B14").Formula = '=""'
$ws.Range("C14").Formula = '=""'
$ws.Range("D14").Formula = '=""'
$ws.Range("E14").Formula = '=""'
$ws.Range("F14").Value = 26634

$ws.Range("A2").Copy()
$ws.Range("A13:A14").PasteSpecial(-4122)
